# DLAD Effective Date January 31, 2022 PROCLTR2022-04
#
# The "30.201-5 Waiver." heading (Heading3 style) needs a bookmark named
# "P30_201_5" wrapped around the "30.201-5 " portion of the text so the
# existing Table-of-Contents hyperlink (which already points at
# w:anchor="P30_201_5") resolves. Word splits the run in two (one run for
# the bookmarked "30.201-5 " text, one for the trailing "Waiver."), and
# automatically renumbers any bookmark ids that collide with the newly
# minted one (the "Part31" bookmark shifts from id 1 to id 2).

$d = $word.ActiveDocument

# Locate the exact "30.201-5 " text (the Heading3 occurrence; the Table of
# Contents entry above it is split into separate "30.201" / "-5" runs via a
# noBreakHyphen, so this literal text only matches the heading we want).
$rng = $d.Content
$found = $rng.Find.Execute("30.201-5 ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Bookmarks.Add("P30_201_5", $rng) | Out-Null
}
